$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptocurrency price/volume update -- Tue Apr 16 15:43:01 UTC 2024

$ws.Range("D2").Value = '62.018.43'
$ws.Range("E2").Value = '  -3.94%  '
$ws.Range("D3").Value = '3.024.24'
$ws.Range("E3").Value = '  -3.97%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").NumberFormat = '@'
$ws.Range("D5").Value = '528.10'
$ws.Range("D5").Style = 'Normal'
$ws.Range("E5").Value = '  -6.33%  '
$ws.Range("D6").NumberFormat = '@'
$ws.Range("D6").Value = '129.74'
$ws.Range("D6").Style = 'Normal'
$ws.Range("E6").Value = '  -9.18%  '
$ws.Range("E7").Value = '  +0.04%  '
$ws.Range("D8").Value = '3.020.81'
$ws.Range("E8").Value = '  -3.84%  '
$ws.Range("E9").Value = '  -1.99%  '
$ws.Range("D10").NumberFormat = '@'
$ws.Range("D10").Value = '0.148'
$ws.Range("D10").Style = 'Normal'
$ws.Range("E10").Value = '  -3.72%  '
$ws.Range("D11").NumberFormat = '@'
$ws.Range("D11").Value = '6.02'
$ws.Range("D11").Style = 'Normal'
$ws.Range("E11").Value = '  -11.57%  '
$ws.Range("D12").NumberFormat = '@'
$ws.Range("D12").Value = '0.444'
$ws.Range("D12").Style = 'Normal'
$ws.Range("E12").Value = '  -4.65%  '
$ws.Range("D13").NumberFormat = '@'
$ws.Range("D13").Value = '0.0000219'
$ws.Range("D13").Style = 'Normal'
$ws.Range("E13").Value = '  -1.28%  '
$ws.Range("D14").NumberFormat = '@'
$ws.Range("D14").Value = '33.27'
$ws.Range("D14").Style = 'Normal'
$ws.Range("E14").Value = '  -9.60%  '
$ws.Range("D15").Value = '3.483.97'
$ws.Range("E15").Value = '  -4.58%  '
$ws.Range("D16").Value = '62.038.10'
$ws.Range("E16").Value = '  -4.00%  '
$ws.Range("E17").Value = '  -2.55%  '
$ws.Range("D18").Value = '3.041.49'
$ws.Range("E18").Value = '  -3.35%  '
$ws.Range("D19").NumberFormat = '@'
$ws.Range("D19").Value = '6.41'
$ws.Range("D19").Style = 'Normal'
$ws.Range("E19").Value = '  -6.25%  '
$ws.Range("D20").NumberFormat = '@'
$ws.Range("D20").Value = '469.94'
$ws.Range("D20").Style = 'Normal'
$ws.Range("E20").Value = '  -9.01%  '
$ws.Range("D21").NumberFormat = '@'
$ws.Range("D21").Value = '12.98'
$ws.Range("D21").Style = 'Normal'
$ws.Range("E21").Value = '  -7.39%  '
$ws.Range("D22").NumberFormat = '@'
$ws.Range("D22").Value = '0.683'
$ws.Range("D22").Style = 'Normal'
$ws.Range("E22").Value = '  -4.53%  '
$ws.Range("D23").NumberFormat = '@'
$ws.Range("D23").Value = '6.95'
$ws.Range("D23").Style = 'Normal'
$ws.Range("E23").Value = '  -6.66%  '
$ws.Range("D24").NumberFormat = '@'
$ws.Range("D24").Value = '77.67'
$ws.Range("D24").Style = 'Normal'
$ws.Range("E24").Value = '  -1.76%  '
$ws.Range("E25").Value = '  -8.38%  '
$ws.Range("D26").NumberFormat = '@'
$ws.Range("D26").Value = '0.998'
$ws.Range("D26").Style = 'Normal'
$ws.Range("E26").Value = '  +0.06%  '
$ws.Range("D27").NumberFormat = '@'
$ws.Range("D27").Value = '2.62'
$ws.Range("D27").Style = 'Normal'
$ws.Range("E27").Value = '  -7.60%  '
$ws.Range("D28").NumberFormat = '@'
$ws.Range("D28").Value = '8.00'
$ws.Range("D28").Style = 'Normal'
$ws.Range("E28").Value = '  -10.13%  '
$ws.Range("E29").Value = '  +0.22%  '
$ws.Range("D30").NumberFormat = '@'
$ws.Range("D30").Value = '25.22'
$ws.Range("D30").Style = 'Normal'
$ws.Range("E30").Value = '  -5.29%  '
$ws.Range("E31").Value = '  -15.98%  '
$ws.Range("E32").Value = '  -4.32%  '
$ws.Range("D33").NumberFormat = '@'
$ws.Range("D33").Value = '2.33'
$ws.Range("D33").Style = 'Normal'
$ws.Range("E33").Value = '  -10.23%  '
$ws.Range("D34").NumberFormat = '@'
$ws.Range("D34").Value = '56.05'
$ws.Range("D34").Style = 'Normal'
$ws.Range("E34").Value = '  +3.88%  '
$ws.Range("D35").NumberFormat = '@'
$ws.Range("D35").Value = '5.10'
$ws.Range("D35").Style = 'Normal'
$ws.Range("E35").Value = '  -5.31%  '
$ws.Range("D36").NumberFormat = '@'
$ws.Range("D36").Value = '5.75'
$ws.Range("D36").Style = 'Normal'
$ws.Range("E36").Value = '  -5.48%  '
$ws.Range("D37").NumberFormat = '@'
$ws.Range("D37").Value = '459.68'
$ws.Range("D37").Style = 'Normal'
$ws.Range("E37").Value = '  -16.55%  '
$ws.Range("D38").Value = '3.048.80'
$ws.Range("E38").Value = '  -3.52%  '
$ws.Range("E39").Value = '  -11.71%  '
$ws.Range("D40").NumberFormat = '@'
$ws.Range("D40").Value = '0.0767'
$ws.Range("D40").Style = 'Normal'
$ws.Range("E40").Value = '  -6.99%  '
$ws.Range("E41").Value = '  -9.03%  '
$ws.Range("D42").NumberFormat = '@'
$ws.Range("D42").Value = '7.87'
$ws.Range("D42").Style = 'Normal'
$ws.Range("E42").Value = '  -5.02%  '
$ws.Range("D43").NumberFormat = '@'
$ws.Range("D43").Value = '2.49'
$ws.Range("D43").Style = 'Normal'
$ws.Range("E43").Value = '  -9.06%  '
$ws.Range("D45").NumberFormat = '@'
$ws.Range("D45").Value = '0.242'
$ws.Range("D45").Style = 'Normal'
$ws.Range("E45").Value = '  -8.65%  '
$ws.Range("E46").Value = '  -11.90%  '
$ws.Range("D47").Value = '0.0₃0509'
$ws.Range("E47").Value = '  -1.72%  '
$ws.Range("D48").NumberFormat = '@'
$ws.Range("D48").Value = '23.73'
$ws.Range("D48").Style = 'Normal'
$ws.Range("E48").Value = '  -6.36%  '
$ws.Range("D49").NumberFormat = '@'
$ws.Range("D49").Value = '0.105'
$ws.Range("D49").Style = 'Normal'
$ws.Range("E49").Value = '  -2.63%  '
$ws.Range("D50").NumberFormat = '@'
$ws.Range("D50").Value = '114.69'
$ws.Range("D50").Style = 'Normal'
$ws.Range("E50").Value = '  -4.89%  '
$ws.Range("B51").Value = 'CoreDAO'
$ws.Range("C51").Value = 'https://coinranking.com/coin/HFvoXUQh4+coredao-core'
$ws.Range("D51").NumberFormat = '@'
$ws.Range("D51").Value = '2.18'
$ws.Range("D51").Style = 'Normal'
$ws.Range("E51").Value = '  -1.66%  '
